$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45867
$ws.Range("B2").Value = 4654.8038038735
$ws.Range("C2").Value = 4572.987492562
$ws.Range("D2").Value = 2376
$ws.Range("E2").Value = 5174.490502
$ws.Range("F2").Value = 113.194757945354

$ws.Range("A3").Value = 45868
$ws.Range("B3").Value = 4654.80380387881
$ws.Range("C3").Value = 4574.55916327919
$ws.Range("D3").Value = 2376
$ws.Range("E3").Value = 5174.490502
$ws.Range("F3").Value = 113.260244225016

$ws.Range("A4").Value = 45869
$ws.Range("B4").Value = 4654.8038038703
$ws.Range("C4").Value = 4539.1800145175
$ws.Range("D4").Value = 2376
$ws.Range("E4").Value = 5174.490502
$ws.Range("F4").Value = 111.786113026967

$ws.Range("A5").Value = 45870
$ws.Range("B5").Value = 5180.74451596793
$ws.Range("C5").Value = 4221.92049285876
$ws.Range("D5").Value = 1944
$ws.Range("E5").Value = 5482.53209
$ws.Range("F5").Value = 107.487836120451

$ws.Range("A6").Value = 45871
$ws.Range("B6").Value = 915.632687220351
$ws.Range("C6").Value = 1586.03949980936
$ws.Range("D6").Value = 1944
$ws.Range("E6").Value = 1753.391387
$ws.Range("F6").Value = 19.9915916495422

$ws.Range("A7").Value = 45872
$ws.Range("B7").Value = 794.873906828036
$ws.Range("C7").Value = 1541.24862209179
$ws.Range("D7").Value = 1944
$ws.Range("E7").Value = 1582.617617
$ws.Range("F7").Value = 16.0413471776565

$ws.Range("A8").Value = 45873
$ws.Range("B8").Value = 3947.0732721826
$ws.Range("C8").Value = 3828.97681251028
$ws.Range("D8").Value = 1944
$ws.Range("E8").Value = 4361.536799
$ws.Range("F8").Value = 95.8100141386534

$ws.Range("A9").Value = 45874
$ws.Range("B9").Value = 3947.0732721826
$ws.Range("C9").Value = 3812.63758512219
$ws.Range("D9").Value = 1944
$ws.Range("E9").Value = 4361.536799
$ws.Range("F9").Value = 95.1292129974828

$ws.Range("A10").Value = 45875
$ws.Range("B10").Value = 3947.0732721826
$ws.Range("C10").Value = 3755.54822468092
$ws.Range("D10").Value = 1944
$ws.Range("E10").Value = 4361.536799
$ws.Range("F10").Value = 92.7504896457631

$ws.Range("A11").Value = 45876
$ws.Range("B11").Value = 3947.0732721826
$ws.Range("C11").Value = 3718.50148601109
$ws.Range("D11").Value = 1944
$ws.Range("E11").Value = 4361.536799
$ws.Range("F11").Value = 91.2068755345204

$ws.Range("A12").Value = 45877
$ws.Range("B12").Value = 3947.0732721826
$ws.Range("C12").Value = 3669.43091028821
$ws.Range("D12").Value = 1944
$ws.Range("E12").Value = 4361.536799
$ws.Range("F12").Value = 89.1622682127335

$ws.Range("A13").Value = 45878
$ws.Range("B13").Value = 719.129912541875
$ws.Range("C13").Value = 1295.01779118882
$ws.Range("D13").Value = 1944
$ws.Range("E13").Value = 1499.56104
$ws.Range("F13").Value = 5.47703827695598

$ws.Range("A14").Value = 45879
$ws.Range("B14").Value = 636.544690493214
$ws.Range("C14").Value = 1237.30835818663
$ws.Range("D14").Value = 1944
$ws.Range("E14").Value = 1408.847898
$ws.Range("F14").Value = 2.73381523722586

$ws.Range("A15").Value = 45880
$ws.Range("B15").Value = 3742.4995935087
$ws.Range("C15").Value = 3509.72202336058
$ws.Range("D15").Value = 1944
$ws.Range("E15").Value = 4164.588359
$ws.Range("F15").Value = 82.8254495354947
